# Swap the full contents of data rows 13 and 14 (the two sightings at
# "Snäckstavik, Srm" got their row order/content exchanged), while keeping
# row 14 beneath row 13 in the sheet.
#
# Helper: write a value into a cell while defeating Excel's "smart" type
# auto-detection for text that merely looks numeric/date/time (e.g. "5",
# "08:18", "2010-01-25"), then strip the leftover quote-prefix formatting
# so no stray style is left behind on the cell.
function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 13 becomes what Row 14 used to be ----
$ws.Range("A13").Value = 130881366
$ws.Range("B13").Value = 57881
Set-TextCell $ws "D13" "NT"
$ws.Range("E13").Value = 100049
Set-TextCell $ws "F13" "Spillkråka"
Set-TextCell $ws "G13" "Dryocopus martius"
Set-TextCell $ws "H13" "(Linnaeus, 1758)"
Set-TextCell $ws "I13" "1"
$ws.Range("J13").ClearContents()
Set-TextCell $ws "K13" "adult"
Set-TextCell $ws "M13" "äldre spår"
$ws.Range("Q13").Value = 656781
$ws.Range("R13").Value = 6559672
Set-TextCell $ws "Z13" "08:48"
Set-TextCell $ws "AB13" "08:48"
Set-TextCell $ws "AX13" "Stuart Fell"

# ---- Row 14 becomes what Row 13 used to be ----
$ws.Range("A14").Value = 130882201
$ws.Range("B14").Value = 93095
Set-TextCell $ws "D14" "LC"
$ws.Range("E14").Value = 4364
Set-TextCell $ws "F14" "Dropptaggsvamp"
Set-TextCell $ws "G14" "Hydnellum ferrugineum"
Set-TextCell $ws "H14" "(Fr.:Fr.) P. Karst."
Set-TextCell $ws "I14" "5"
Set-TextCell $ws "J14" "fruktkroppar"
$ws.Range("K14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("Q14").Value = 656955
$ws.Range("R14").Value = 6559350
Set-TextCell $ws "Z14" "08:18"
Set-TextCell $ws "AB14" "08:18"
Set-TextCell $ws "AX14" "Stuart Fell, Liam Martin"
